# Add a new row (41) to the "Bereiche" lookup table, matching the
# formatting of the last existing data row (40).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 40 (A:D) down into the new row 41 before
# writing values, so the new row picks up the same style (border/fill/
# font/wrap) as the rest of the table instead of the workbook default.
$ws.Range("A40:D40").Copy()
$ws.Range("A41:D41").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A41").Value = "Z03_B04"
$ws.Range("B41").Value = "Z03"
$ws.Range("C41").Value = "Soziale Lage und Gesundheit"
$ws.Range("D41").Value = "X"
